$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Update actual start/duration and percent complete for rows 22-33
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1

$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1

$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 0.75

$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 1

$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.5

$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1

$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1

$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 0.25

$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 0

$ws.Range("G32").Value = 0
$ws.Range("G33").Value = 0

$ws.Range("G35").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 0

# View changes
$ws.Activate()
$ws.Range("G33").Select()

# Remove Sheet1
$excel.DisplayAlerts = $false
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Delete()
